# Conference-call bingo workbook edit
# Adds quotation marks to phrase entries (vs. bare jargon/sound-effect rows)
# on the "list" sheet, reshuffles the word order, adds two new phrases, and
# removes a few retired ones -- matching the upstream "BINGO_cc.xlsx" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

$ws.Cells.Item(1, 1).Value = 'Quotes'
$ws.Cells.Item(2, 1).Value = '"Happy [weekday]!"'
$ws.Cells.Item(3, 1).Value = '"My outlook / WebEx is not working"'
$ws.Cells.Item(4, 1).Value = '"I have to jump to another call"'
$ws.Cells.Item(5, 1).Value = '"[We] can''t see your screen"'
$ws.Cells.Item(6, 1).Value = 'typing noise'
$ws.Cells.Item(7, 1).Value = '"Can you email that to everyone?"'
$ws.Cells.Item(8, 1).Value = '"Sorry I was on Mute"'
$ws.Cells.Item(9, 1).Value = '"It''s [almost] Friday!"'
$ws.Cells.Item(10, 1).Value = 'echo / feedback'
$ws.Cells.Item(11, 1).Value = '"[Let''s] take this offline"'
$ws.Cells.Item(12, 1).Value = '"I have a hard stop"'
$ws.Cells.Item(13, 1).Value = '"Who just joined?"'
$ws.Cells.Item(14, 1).Value = '"Go ahead" (talking at once)'
$ws.Cells.Item(15, 1).Value = '"Is ____ on the call?"'
$ws.Cells.Item(16, 1).Value = '"Your phone was breaking up"'
$ws.Cells.Item(17, 1).Value = '"Can you repeat the question?"'
$ws.Cells.Item(18, 1).Value = 'wind / road noise'
$ws.Cells.Item(19, 1).Value = '"Please mute if you aren''t talking"'
$ws.Cells.Item(20, 1).Value = '"Can you hear me?"'
$ws.Cells.Item(21, 1).Value = '"Can you / everyone see my screen?"'
$ws.Cells.Item(22, 1).Value = '"You''re not (are you) sharing(?)"'
$ws.Cells.Item(23, 1).Value = '"Let me IM them to see if they''re joining"'
$ws.Cells.Item(24, 1).Value = '"I''ll have to get back to you"'
$ws.Cells.Item(25, 1).Value = '"Ok. Let''s get started"'
$ws.Cells.Item(26, 1).Value = '"I''ll give you back n minutes "'
$ws.Cells.Item(27, 1).Value = '"Reach out to  ____"'
$ws.Cells.Item(28, 1).Value = '(Starting) "This won''t take the whole time"'
$ws.Cells.Item(29, 1).Value = '"level-set"'
$ws.Cells.Item(30, 1).Value = '"Sorry I’m late" (lame excuse)'
$ws.Cells.Item(31, 1).Value = '"I was multi-tasking"'
$ws.Cells.Item(32, 1).Value = '"low hanging fruit"'
$ws.Cells.Item(33, 1).Value = '"at the end of the day"'
$ws.Cells.Item(34, 1).Value = '"win-win"'
$ws.Cells.Item(35, 1).Value = '"Loop in  ____"'
$ws.Cells.Item(36, 1).Value = '"Keep ____ in the loop"'
$ws.Cells.Item(37, 1).Value = 'awkward silence'
$ws.Cells.Item(38, 1).Value = '"circle back"'
$ws.Cells.Item(39, 1).Value = '"touch base"'
$ws.Cells.Item(40, 1).Value = '"I’ll take silence as ____"'
$ws.Cells.Item(41, 1).Value = '"You''re breaking up"'
$ws.Cells.Item(42, 1).Value = '"I''m having computer problems"'
$ws.Cells.Item(43, 1).Value = '"on the same page"'
$ws.Cells.Item(44, 1).Value = '"I''m having technical difficulties"'
$ws.Cells.Item(45, 1).Value = '"I have to drop"'
$ws.Cells.Item(46, 1).Value = 'unnecessary verbing ("Let''s solution that")'
$ws.Cells.Item(47, 1).Value = '"on/off the radar"'
$ws.Cells.Item(48, 1).Value = '"It is what it is"'
$ws.Cells.Item(49, 1).Value = '"ping me / you"'
$ws.Cells.Item(50, 1).Value = '"have the bandwidth" (i.e. work capacity)'
$ws.Cells.Item(51, 1).Value = '"ducks in a row"'
$ws.Cells.Item(52, 1).Value = 'chewing sounds'
$ws.Cells.Item(53, 1).Value = 'animal sounds'
$ws.Cells.Item(54, 1).Value = '"You''re cutting out"'
$ws.Cells.Item(55, 1).Value = '"It''s loading"'
$ws.Cells.Item(56, 1).Value = '"Next slide please"'
$ws.Cells.Item(57, 1).Value = 'child sounds'
$ws.Cells.Item(58, 1).Value = '"voluntold"'
$ws.Cells.Item(59, 1).Value = '(video) never looking at camera'
$ws.Cells.Item(60, 1).Value = '(video) animal(s) walking around'
$ws.Cells.Item(61, 1).Value = '(video) partner in view'
$ws.Cells.Item(62, 1).Value = '(video) kid(s) in view/interrupting'

# Restore the selection/scroll state recorded in the authored workbook.
$ws.Activate()
$ws.Range("A59").Select()
$excel.ActiveWindow.ScrollRow = 34

# Note: xl/worksheets/sheet1.xml B2 holds a volatile formula
#   ="Happy "&TEXT(TODAY(),"dddd")
# its cached value tracks the live clock; left untouched so recalculation
# (which the harness drives with the capture-time clock) reproduces the
# authored "Happy Thursday" value naturally.

# Note: the "mon1"/"grid" sheets (B1/F3/E4/D5) and the reordering of
# xl/sharedStrings.xml are pure internal shared-string reindexing side
# effects of the above edits; the displayed text in those cells is
# unchanged, so no further action is required there.
